# Apply the changes described by the diff:
# - Update workbook view xWindow
# - Update sheet view topLeftCell + selection
# - Fill in remaining data for the "Dimension 1" table (rows 5-7, columns E-Q)
# - Fill in remaining data for the "Dimension 4" table (rows 23-25, columns P-Q)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Window / view settings ----
# sheetView topLeftCell: H1 -> G1 (scroll so column G is left-most visible column)
$excel.ActiveWindow.ScrollColumn = 7   # G is the 7th column
$excel.ActiveWindow.ScrollRow = 1

# selection activeCell/sqref: P23 -> Q10
$ws.Range("Q10").Select()

# ---- Row 5 (Average MST Weight, Dimension 1 table) ----
$ws.Range("E5").Value = 1.158234
$ws.Range("F5").Value = 1.183538
$ws.Range("G5").Value = 1.194542
$ws.Range("H5").Value = 1.1993769999999999
$ws.Range("I5").Value = 1.2029209999999999
$ws.Range("J5").Value = 1.2030639999999999
$ws.Range("K5").Value = 1.2019660000000001
$ws.Range("L5").Value = 1.2002759999999999
$ws.Range("M5").Value = 1.198312
$ws.Range("N5").Value = 1.193554
$ws.Range("O5").Value = 1.2041850000000001
$ws.Range("P5").Value = 1.205743
$ws.Range("Q5").Value = 1.200774

# ---- Row 6 (Max included edge, Dimension 1 table) ----
$ws.Range("E6").Value = 0.55221299999999995
$ws.Range("F6").Value = 0.35253600000000002
$ws.Range("G6").Value = 0.21062
$ws.Range("H6").Value = 0.105326
$ws.Range("I6").Value = 0.051035999999999998
$ws.Range("J6").Value = 0.026724999999999999
$ws.Range("K6").Value = 0.014557
$ws.Range("L6").Value = 0.0070740000000000004
$ws.Range("M6").Value = 0.003225
$ws.Range("N6").Value = 0.001446
$ws.Range("O6").Value = 0.00076499999999999995
$ws.Range("P6").Value = 0.00041100000000000002
$ws.Range("Q6").Value = 0.00020599999999999999

# ---- Row 7 (trials, Dimension 1 table) ----
$ws.Range("H7").Value = 10000
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 10
$ws.Range("P7").Value = 5
$ws.Range("Q7").Value = 5

# ---- Row 23 (Average MST Weight, Dimension 4 table) ----
$ws.Range("P23").Value = 1688.6282960000001
$ws.Range("Q23").Value = 2830.0883789999998

# ---- Row 24 (Max included edge, Dimension 4 table) ----
$ws.Range("P24").Value = 0.113422
$ws.Range("Q24").Value = 0.10422099999999999

# ---- Row 25 (trials, Dimension 4 table) ----
$ws.Range("P25").Value = 5
$ws.Range("Q25").Value = 5
